# Revert "added slide and slight text change pmp"
#
# 1. Remove the last slide (the "Vragen?" / Questions slide that was added).
# 2. Restore the cached Date and Slide Number placeholder text on the
#    slide master and every slide layout (20-1-2021 -> 18-1-2021,
#    ‹#› -> ‹nr.›).

$p = $ppt.ActivePresentation

# --- 1. Delete the last slide (slide 8 / "Vragen?") -----------------------
$lastIndex = $p.Slides.Count
if ($lastIndex -gt 0) {
    $p.Slides.Item($lastIndex).Delete()
}

# --- 2. Restore footer date / slide-number placeholder text ---------------
$oldDate = "20-1-2021"
$newDate = "18-1-2021"
$oldNum  = [char]0x2039 + "#" + [char]0x203A
$newNum  = [char]0x2039 + "nr." + [char]0x203A

function Update-FooterFields($container) {
    $shapes = $container.Shapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if (-not $sh.HasTextFrame) { continue }
        if ($sh.Type -ne 14) { continue }
        $phType = $sh.PlaceholderFormat.Type
        $tr = $sh.TextFrame.TextRange
        if ($phType -eq 16) {
            # ppPlaceholderDate
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        } elseif ($phType -eq 13) {
            # ppPlaceholderSlideNumber
            if ($tr.Text -eq $oldNum) {
                $tr.Text = $newNum
            }
        }
    }
}

# Slide master
Update-FooterFields($p.SlideMaster)

# All slide layouts
$layouts = $p.SlideMaster.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    Update-FooterFields($layouts.Item($L))
}
